$d = $word.ActiveDocument
$d.Content.Find.Execute("dsa", $true, $false, $false, $false, $false,
                         $true, 1, $false, "[PLACEHOLDER]", 2)
